# Daily attendance processing - swap order of first two "Recorded By"
# contributors in column G (keeps any trailing extra entries, e.g. "system",
# in their original position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    $val = [string]$val
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ",\s*"
    if ($parts.Count -ge 2) {
        $swapped = @($parts[1], $parts[0])
        if ($parts.Count -gt 2) {
            $swapped += $parts[2..($parts.Count - 1)]
        }
        $cell.Value = [string]::Join(", ", $swapped)
    }
}
